$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Intel(R) Wi-Fi 6E AX211 160MHz - 22.190.0.4
$ws.Range("C3").Value = 2981
$ws.Range("D3").Value = 94.5

# Row 4 - Intel(R) Wi-Fi 6E AX211 160MHz - 23.40.0.4
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 492
$ws.Range("D4").Value = 96.8

# Row 5 - Totals
$ws.Range("B5").Value = 11
$ws.Range("C5").Value = 3473

# Row 13 - clear Driver Vintage date value
$ws.Range("E13").Value = $null

# Row 15 - Total Samples update
$ws.Range("B15").Value = 265400
